$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 4 = "Thank you for your attention!" text box holding the
# Github-workshop instructions on slide 1.
$shape = $s.Shapes.Item(4)
$tr = $shape.TextFrame.TextRange

# --- Paragraph 2: fix the repository URL -----------------------------
$para2 = $tr.Paragraphs(2)
$urlRun = $para2.Runs(2)
$urlRun.Text = "https://github.com/IEECR/Git-Github-Workshop/tree/main"

# --- Paragraph 3: split the instructions sentence into several runs --
# so the folder/file names can be italicised.
$para3 = $tr.Paragraphs(3)
$run1 = $para3.Runs(1)
$run1.Text = "In the repository move to the folder "

$null = $para3.InsertAfter("sessions_1")
$null = $para3.InsertAfter(" > ")
$null = $para3.InsertAfter("instructions")
$null = $para3.InsertAfter(" and click on ")
$null = $para3.InsertAfter("module_1.md")

# Italicise the folder / file name runs (2, 4 and 6) only - do this
# after all the runs have been created so the italic flag does not leak
# onto the plain-text runs in between.
$para3.Runs(2).Font.Italic = $true
$para3.Runs(4).Font.Italic = $true
$para3.Runs(6).Font.Italic = $true
